$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4047206666666667
$ws.Range("H2").Value = 1.214162
$ws.Range("I2").Value = 0.02578034976888792
$ws.Range("J2").Value = 0.02578034976888792
$ws.Range("M2").Value = 4.861952666666666
$ws.Range("N2").Value = 14.585858
$ws.Range("O2").Value = 0.3995648519435639
$ws.Range("P2").Value = 0.3995648519435638
$ws.Range("Q2").Value = 1.967732724555111
$ws.Range("R2").Value = 17.709594520996
$ws.Range("S2").Value = 0.01030092163845899
$ws.Range("T2").Value = 0.01030092163845899
$ws.Range("G3").Value = 0.4047206666666667
$ws.Range("H3").Value = 1.214162
$ws.Range("I3").Value = 0.02578034976888792
$ws.Range("J3").Value = 0.02578034976888792
$ws.Range("O3").Value = 0.04932556406896855
$ws.Range("P3").Value = 0.04932556406896854
$ws.Range("Q3").Value = 0.2429130743195556
$ws.Range("R3").Value = 2.186217668876
$ws.Range("S3").Value = 0.0012716302942457
$ws.Range("T3").Value = 0.001271630294245699
$ws.Range("G4").Value = 0.4047206666666667
$ws.Range("H4").Value = 1.214162
$ws.Range("I4").Value = 0.02578034976888792
$ws.Range("J4").Value = 0.02578034976888792
$ws.Range("M4").Value = 4.206754333333333
$ws.Range("N4").Value = 12.620263
$ws.Range("O4").Value = 0.3457193616641432
$ws.Range("P4").Value = 0.3457193616641432
$ws.Range("Q4").Value = 1.702560418289555
$ws.Range("R4").Value = 15.323043764606
$ws.Range("S4").Value = 0.008912766065578274
$ws.Range("T4").Value = 0.008912766065578274
$ws.Range("G5").Value = 0.4047206666666667
$ws.Range("H5").Value = 1.214162
$ws.Range("I5").Value = 0.02578034976888792
$ws.Range("J5").Value = 0.02578034976888792
$ws.Range("M5").Value = 2.499212666666667
$ws.Range("N5").Value = 7.497638
$ws.Range("O5").Value = 0.2053902223233243
$ws.Range("P5").Value = 0.2053902223233243
$ws.Range("Q5").Value = 1.011483016595111
$ws.Range("R5").Value = 9.103347149356
$ws.Range("S5").Value = 0.005295031770604953
$ws.Range("T5").Value = 0.005295031770604953
$ws.Range("I6").Value = 0.03222381288358415
$ws.Range("J6").Value = 0.03222381288358415
$ws.Range("M6").Value = 4.861952666666666
$ws.Range("N6").Value = 14.585858
$ws.Range("O6").Value = 0.3995648519435639
$ws.Range("P6").Value = 0.3995648519435638
$ws.Range("Q6").Value = 2.459541925900889
$ws.Range("R6").Value = 22.135877333108
$ws.Range("S6").Value = 0.01287550302388641
$ws.Range("T6").Value = 0.01287550302388641
$ws.Range("I7").Value = 0.03222381288358415
$ws.Range("J7").Value = 0.03222381288358415
$ws.Range("O7").Value = 0.04932556406896855
$ws.Range("P7").Value = 0.04932556406896854
$ws.Range("S7").Value = 0.001589457746935684
$ws.Range("T7").Value = 0.001589457746935684
$ws.Range("I8").Value = 0.03222381288358415
$ws.Range("J8").Value = 0.03222381288358415
$ws.Range("M8").Value = 4.206754333333333
$ws.Range("N8").Value = 12.620263
$ws.Range("O8").Value = 0.3457193616641432
$ws.Range("P8").Value = 0.3457193616641432
$ws.Range("Q8").Value = 2.128093250626444
$ws.Range("R8").Value = 19.152839255638
$ws.Range("S8").Value = 0.01114039602049751
$ws.Range("T8").Value = 0.01114039602049751
$ws.Range("I9").Value = 0.03222381288358415
$ws.Range("J9").Value = 0.03222381288358415
$ws.Range("M9").Value = 2.499212666666667
$ws.Range("N9").Value = 7.497638
$ws.Range("O9").Value = 0.2053902223233243
$ws.Range("P9").Value = 0.2053902223233243
$ws.Range("Q9").Value = 1.264290040820889
$ws.Range("R9").Value = 11.378610367388
$ws.Range("S9").Value = 0.006618456092264551
$ws.Range("T9").Value = 0.006618456092264552
$ws.Range("G10").Value = 1.039987
$ws.Range("H10").Value = 3.119961
$ws.Range("I10").Value = 0.06624625531460326
$ws.Range("J10").Value = 0.06624625531460326
$ws.Range("M10").Value = 4.861952666666666
$ws.Range("N10").Value = 14.585858
$ws.Range("O10").Value = 0.3995648519435639
$ws.Range("P10").Value = 0.3995648519435638
$ws.Range("Q10").Value = 5.056367567948667
$ws.Range("R10").Value = 45.50730811153799
$ws.Range("S10").Value = 0.02646967519659498
$ws.Range("T10").Value = 0.02646967519659498
$ws.Range("G11").Value = 1.039987
$ws.Range("H11").Value = 3.119961
$ws.Range("I11").Value = 0.06624625531460326
$ws.Range("J11").Value = 0.06624625531460326
$ws.Range("O11").Value = 0.04932556406896855
$ws.Range("P11").Value = 0.04932556406896854
$ws.Range("Q11").Value = 0.6241995040753334
$ws.Range("R11").Value = 5.617795536678001
$ws.Range("S11").Value = 0.003267633910849711
$ws.Range("T11").Value = 0.003267633910849711
$ws.Range("G12").Value = 1.039987
$ws.Range("H12").Value = 3.119961
$ws.Range("I12").Value = 0.06624625531460326
$ws.Range("J12").Value = 0.06624625531460326
$ws.Range("M12").Value = 4.206754333333333
$ws.Range("N12").Value = 12.620263
$ws.Range("O12").Value = 0.3457193616641432
$ws.Range("P12").Value = 0.3457193616641432
$ws.Range("Q12").Value = 4.374969818860333
$ws.Range("R12").Value = 39.374728369743
$ws.Range("S12").Value = 0.0229026131000045
$ws.Range("T12").Value = 0.0229026131000045
$ws.Range("G13").Value = 1.039987
$ws.Range("H13").Value = 3.119961
$ws.Range("I13").Value = 0.06624625531460326
$ws.Range("J13").Value = 0.06624625531460326
$ws.Range("M13").Value = 2.499212666666667
$ws.Range("N13").Value = 7.497638
$ws.Range("O13").Value = 0.2053902223233243
$ws.Range("P13").Value = 0.2053902223233243
$ws.Range("Q13").Value = 2.599148683568667
$ws.Range("R13").Value = 23.392338152118
$ws.Range("S13").Value = 0.01360633310715407
$ws.Range("T13").Value = 0.01360633310715407
$ws.Range("G14").Value = 13.74822133333333
$ws.Range("H14").Value = 41.244664
$ws.Range("I14").Value = 0.8757495820329246
$ws.Range("J14").Value = 0.8757495820329247
$ws.Range("M14").Value = 4.861952666666666
$ws.Range("N14").Value = 14.585858
$ws.Range("O14").Value = 0.3995648519435639
$ws.Range("P14").Value = 0.3995648519435638
$ws.Range("Q14").Value = 66.84320137352356
$ws.Range("R14").Value = 601.5888123617119
$ws.Range("S14").Value = 0.3499187520846235
$ws.Range("T14").Value = 0.3499187520846235
$ws.Range("G15").Value = 13.74822133333333
$ws.Range("H15").Value = 41.244664
$ws.Range("I15").Value = 0.8757495820329246
$ws.Range("J15").Value = 0.8757495820329247
$ws.Range("O15").Value = 0.04932556406896855
$ws.Range("P15").Value = 0.04932556406896854
$ws.Range("Q15").Value = 8.251673278785779
$ws.Range("R15").Value = 74.265059509072
$ws.Range("S15").Value = 0.04319684211693745
$ws.Range("T15").Value = 0.04319684211693745
$ws.Range("G16").Value = 13.74822133333333
$ws.Range("H16").Value = 41.244664
$ws.Range("I16").Value = 0.8757495820329246
$ws.Range("J16").Value = 0.8757495820329247
$ws.Range("M16").Value = 4.206754333333333
$ws.Range("N16").Value = 12.620263
$ws.Range("O16").Value = 0.3457193616641432
$ws.Range("P16").Value = 0.3457193616641432
$ws.Range("Q16").Value = 57.83538966962577
$ws.Range("R16").Value = 520.5185070266319
$ws.Range("S16").Value = 0.302763586478063
$ws.Range("T16").Value = 0.302763586478063
$ws.Range("G17").Value = 13.74822133333333
$ws.Range("H17").Value = 41.244664
$ws.Range("I17").Value = 0.8757495820329246
$ws.Range("J17").Value = 0.8757495820329247
$ws.Range("M17").Value = 2.499212666666667
$ws.Range("N17").Value = 7.497638
$ws.Range("O17").Value = 0.2053902223233243
$ws.Range("P17").Value = 0.2053902223233243
$ws.Range("Q17").Value = 34.35972890040356
$ws.Range("R17").Value = 309.237560103632
$ws.Range("S17").Value = 0.1798704013533007
$ws.Range("T17").Value = 0.1798704013533008
